# Add a new "Compact List" paragraph style (w:styleId="CompactList"),
# mirroring the existing "Compact" style: based on Body Text, marked as
# a quick style, with the same tight before/after spacing (36 twips = 1.8pt).
$d = $word.ActiveDocument

$style = $d.Styles.Add("Compact List", 1)
$style.BaseStyle = "BodyText"
$style.QuickStyle = $true

$pf = $style.ParagraphFormat
$pf.SpaceBefore = 1.8
$pf.SpaceAfter = 1.8
